$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 23 (pushes existing rows 23-27 down to 24-28),
# inheriting the number format of the row being displaced.
$ws.Rows.Item(23).Insert()

# Fill in the new CMIP3-routine GWL data row (GFDL-CM3, rcp85)
$ws.Cells.Item(23, 1).Value = "GFDL-CM3_rcp85_r1i1p1_200601-210012"
$ws.Cells.Item(23, 2).Value = 2023
$ws.Cells.Item(23, 3).Value = 2035
$ws.Cells.Item(23, 4).Value = 2055
$ws.Cells.Item(23, 5).Value = 2071
$ws.Cells.Item(23, 6).Value = "2013-2032"
$ws.Cells.Item(23, 7).Value = "2025-2044"
$ws.Cells.Item(23, 8).Value = "2045-2064"
$ws.Cells.Item(23, 9).Value = "2061-2080"
$ws.Cells.Item(23, 10).Value = 1.2547776910347199
$ws.Cells.Item(23, 11).Value = 1.8814544179047701
$ws.Cells.Item(23, 12).Value = 3.3167032989732199
$ws.Cells.Item(23, 13).Value = 4.9136741564842401

# New routines area further down the sheet: a cell formatted with higher
# numeric precision for upcoming CMIP3 computations.
$ws.Cells.Item(49, 10).NumberFormat = "0.000000000"

# Narrow the J:M columns (now holding shorter GWL values) and update the
# last-saved selection.
$ws.Range("J1:M1").EntireColumn.ColumnWidth = 5.6640625
[void]$ws.Range("Q10").Select()
